$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 9297.637000000001
$ws.Range("I6").Value = 10197.9
$ws.Range("K6").Value = 30593.7
$ws.Range("M6").Value = -30481.7
$ws.Range("H18").Value = 361.8
$ws.Range("I18").Value = 403
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 403
$ws.Range("L18").Value = 300
$ws.Range("M18").Value = -119
$ws.Range("N18").Value = -868
$ws.Range("H19").Value = 2445.8235
$ws.Range("I19").Value = 2118.2727
$ws.Range("K19").Value = 2118.2727
$ws.Range("M19").Value = -1943.2727
$ws.Range("H40").Value = 21014.79
$ws.Range("I40").Value = 19029.615
$ws.Range("J40").Value = 25316
$ws.Range("K40").Value = 19029.615
$ws.Range("L40").Value = 25316
$ws.Range("M40").Value = -18854.615
$ws.Range("N40").Value = -25666
$ws.Range("H98").Value = 1754.5385
$ws.Range("I98").Value = 1754.5385
$ws.Range("K98").Value = 1754.5385
$ws.Range("M98").Value = -256.5385000000001
$ws.Range("H100").Value = 2847.889
$ws.Range("I100").Value = 2252.6365
$ws.Range("J100").Value = 3783.2856
$ws.Range("K100").Value = 2252.6365
$ws.Range("L100").Value = 3783.2856
$ws.Range("M100").Value = -1711.6365
$ws.Range("N100").Value = -4865.2856
$ws.Range("H101").Value = 2836.4285
$ws.Range("J101").Value = 5355.143
$ws.Range("L101").Value = 16065.429
$ws.Range("N101").Value = -19309.429
$ws.Range("H112").Value = 2284.1924
$ws.Range("J112").Value = 2528.0952
$ws.Range("L112").Value = 7584.285600000001
$ws.Range("N112").Value = -9800.285600000001
$ws.Range("H116").Value = 16675460
$ws.Range("I116").Value = 20841596
$ws.Range("K116").Value = 20841596
$ws.Range("M116").Value = -20838154
$ws.Range("H122").Value = 1754.5385
$ws.Range("I122").Value = 1754.5385
$ws.Range("K122").Value = 5263.6155
$ws.Range("M122").Value = -2813.6155
$ws.Range("H132").Value = 177787.31
$ws.Range("I132").Value = 280835.97
$ws.Range("J132").Value = 23214.318
$ws.Range("K132").Value = 842507.9099999999
$ws.Range("L132").Value = 69642.954
$ws.Range("M132").Value = -839977.9099999999
$ws.Range("N132").Value = -74702.954
$ws.Range("H135").Value = 2641.0293
$ws.Range("I135").Value = 1109.6666
$ws.Range("K135").Value = 9986.999400000001
$ws.Range("M135").Value = -7451.999400000001
$ws.Range("H137").Value = 360449.4
$ws.Range("I137").Value = 479404
$ws.Range("J137").Value = 3585.5715
$ws.Range("K137").Value = 1438212
$ws.Range("L137").Value = 10756.7145
$ws.Range("M137").Value = -1435662
$ws.Range("N137").Value = -15856.7145
$ws.Range("H138").Value = 2907.1216
$ws.Range("I138").Value = 913.7273
$ws.Range("J138").Value = 5830.7666
$ws.Range("K138").Value = 2741.1819
$ws.Range("L138").Value = 17492.2998
$ws.Range("M138").Value = 2398.8181
$ws.Range("N138").Value = -27772.2998
$ws.Range("H141").Value = 3220.8684
$ws.Range("I141").Value = 3248.8
$ws.Range("J141").Value = 3116.125
$ws.Range("K141").Value = 9746.400000000001
$ws.Range("L141").Value = 9348.375
$ws.Range("M141").Value = -4566.400000000001
$ws.Range("N141").Value = -19708.375

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 72727
$ws.Range("I2").Value = 72727
$ws.Range("K2").Value = 72727
$ws.Range("M2").Value = -72614
$ws.Range("H32").Value = 11641
$ws.Range("I32").Value = 11531.486
$ws.Range("J32").Value = 12204.214
$ws.Range("K32").Value = 11531.486
$ws.Range("L32").Value = 12204.214
$ws.Range("M32").Value = -11244.486
$ws.Range("N32").Value = -12778.214
$ws.Range("H45").Value = 2843
$ws.Range("I45").Value = 1927.375
$ws.Range("K45").Value = 1927.375
$ws.Range("M45").Value = -1550.375
$ws.Range("H61").Value = 3084.4187
$ws.Range("I61").Value = 1961.5151
$ws.Range("K61").Value = 1961.5151
$ws.Range("M61").Value = -1749.5151
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 1678.7273
$ws.Range("I74").Value = 919.5
$ws.Range("K74").Value = 919.5
$ws.Range("M74").Value = -45.5
$ws.Range("H77").Value = 1678.7273
$ws.Range("I77").Value = 919.5
$ws.Range("K77").Value = 4597.5
$ws.Range("M77").Value = -229.5
$ws.Range("H102").Value = 1323.5834
$ws.Range("J102").Value = 750
$ws.Range("L102").Value = 750
$ws.Range("N102").Value = -3994
$ws.Range("H110").Value = 12497.5
$ws.Range("I110").Value = 8996.25
$ws.Range("K110").Value = 8996.25
$ws.Range("M110").Value = -6951.25
$ws.Range("H116").Value = 72727
$ws.Range("I116").Value = 72727
$ws.Range("K116").Value = 72727
$ws.Range("M116").Value = -70433
$ws.Range("H122").Value = 3721.0217
$ws.Range("I122").Value = 2566.2974
$ws.Range("J122").Value = 8468.223
$ws.Range("K122").Value = 7698.8922
$ws.Range("L122").Value = 25404.669
$ws.Range("M122").Value = -5248.8922
$ws.Range("N122").Value = -30304.669
$ws.Range("H132").Value = 14602.479
$ws.Range("I132").Value = 16919.861
$ws.Range("J132").Value = 6259.9
$ws.Range("K132").Value = 50759.583
$ws.Range("L132").Value = 18779.7
$ws.Range("M132").Value = -48229.583
$ws.Range("N132").Value = -23839.7
$ws.Range("H136").Value = 3084.4187
$ws.Range("I136").Value = 1961.5151
$ws.Range("K136").Value = 5884.5453
$ws.Range("M136").Value = -3334.5453

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 72727
$ws.Range("I3").Value = 72727
$ws.Range("K3").Value = 72727
$ws.Range("M3").Value = -72613
$ws.Range("H20").Value = 2328.4482
$ws.Range("I20").Value = 2491.7273
$ws.Range("J20").Value = 1815.2858
$ws.Range("K20").Value = 2491.7273
$ws.Range("L20").Value = 1815.2858
$ws.Range("M20").Value = -2244.7273
$ws.Range("N20").Value = -2309.2858
$ws.Range("H22").Value = 999
$ws.Range("I22").Value = 999
$ws.Range("K22").Value = 999
$ws.Range("M22").Value = -826
$ws.Range("H86").Value = 2749.5
$ws.Range("I86").Value = 2833
$ws.Range("K86").Value = 2833
$ws.Range("M86").Value = -1710
$ws.Range("H89").Value = 2749.5
$ws.Range("I89").Value = 2833
$ws.Range("K89").Value = 14165
$ws.Range("M89").Value = -8549
$ws.Range("H94").Value = 2518.2856
$ws.Range("I94").Value = 2271.3333
$ws.Range("K94").Value = 2271.3333
$ws.Range("M94").Value = -1820.3333
$ws.Range("H99").Value = 1840.1177
$ws.Range("I99").Value = 1652.5385
$ws.Range("K99").Value = 1652.5385
$ws.Range("M99").Value = -154.5385000000001
$ws.Range("H105").Value = 3649.5715
$ws.Range("I105").Value = 3105.2
$ws.Range("K105").Value = 3105.2
$ws.Range("M105").Value = -1358.2
$ws.Range("H107").Value = 2847.6775
$ws.Range("I107").Value = 2120.4583
$ws.Range("J107").Value = 5341
$ws.Range("K107").Value = 2120.4583
$ws.Range("L107").Value = 5341
$ws.Range("M107").Value = -200.4582999999998
$ws.Range("N107").Value = -9181
$ws.Range("H134").Value = 2584.2563
$ws.Range("I134").Value = 1792.7667
$ws.Range("J134").Value = 5222.5557
$ws.Range("K134").Value = 5378.300099999999
$ws.Range("L134").Value = 15667.6671
$ws.Range("M134").Value = -2843.300099999999
$ws.Range("N134").Value = -20737.6671

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3732
$ws.Range("I16").Value = 3102.75
$ws.Range("J16").Value = 4571
$ws.Range("K16").Value = 3102.75
$ws.Range("L16").Value = 4571
$ws.Range("M16").Value = -2815.75
$ws.Range("N16").Value = -5145
$ws.Range("H31").Value = 4696.6787
$ws.Range("I31").Value = 2969.7334
$ws.Range("J31").Value = 6689.3076
$ws.Range("K31").Value = 2969.7334
$ws.Range("L31").Value = 6689.3076
$ws.Range("M31").Value = -2674.7334
$ws.Range("N31").Value = -7279.3076
$ws.Range("H34").Value = 4696.6787
$ws.Range("I34").Value = 2969.7334
$ws.Range("J34").Value = 6689.3076
$ws.Range("K34").Value = 2969.7334
$ws.Range("L34").Value = 6689.3076
$ws.Range("M34").Value = -2767.7334
$ws.Range("N34").Value = -7093.3076
$ws.Range("H58").Value = 1002028.9
$ws.Range("I58").Value = 1896.7142
$ws.Range("K58").Value = 1896.7142
$ws.Range("M58").Value = -1693.7142
$ws.Range("H86").Value = 13216.968
$ws.Range("I86").Value = 12815.947
$ws.Range("K86").Value = 12815.947
$ws.Range("M86").Value = -11692.947
$ws.Range("H89").Value = 13216.968
$ws.Range("I89").Value = 12815.947
$ws.Range("K89").Value = 64079.735
$ws.Range("M89").Value = -58463.735
$ws.Range("H99").Value = 10171.946
$ws.Range("I99").Value = 10888.389
$ws.Range("J99").Value = 9493.210999999999
$ws.Range("K99").Value = 10888.389
$ws.Range("L99").Value = 9493.210999999999
$ws.Range("M99").Value = -9390.388999999999
$ws.Range("N99").Value = -12489.211
$ws.Range("H113").Value = 3732
$ws.Range("I113").Value = 3102.75
$ws.Range("J113").Value = 4571
$ws.Range("K113").Value = 3102.75
$ws.Range("L113").Value = 4571
$ws.Range("M113").Value = -932.75
$ws.Range("N113").Value = -8911
$ws.Range("H122").Value = 4054.5789
$ws.Range("I122").Value = 1925.9231
$ws.Range("K122").Value = 5777.7693
$ws.Range("M122").Value = -3327.7693
$ws.Range("H126").Value = 10171.946
$ws.Range("I126").Value = 10888.389
$ws.Range("J126").Value = 9493.210999999999
$ws.Range("K126").Value = 32665.167
$ws.Range("L126").Value = 28479.633
$ws.Range("M126").Value = -30195.167
$ws.Range("N126").Value = -33419.633
$ws.Range("H132").Value = 9538431
$ws.Range("I132").Value = 11503899
$ws.Range("J132").Value = 38673
$ws.Range("K132").Value = 34511697
$ws.Range("L132").Value = 116019
$ws.Range("M132").Value = -34509167
$ws.Range("N132").Value = -121079
$ws.Range("H134").Value = 2045.7273
$ws.Range("I134").Value = 1995.7428
$ws.Range("J134").Value = 2240.111
$ws.Range("K134").Value = 5987.2284
$ws.Range("L134").Value = 6720.333
$ws.Range("M134").Value = -3452.2284
$ws.Range("N134").Value = -11790.333
$ws.Range("H136").Value = 1002028.9
$ws.Range("I136").Value = 1896.7142
$ws.Range("K136").Value = 5690.142599999999
$ws.Range("M136").Value = -3140.142599999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1596.1818
$ws.Range("J2").Value = 3287.6
$ws.Range("L2").Value = 19725.6
$ws.Range("N2").Value = -19951.6
$ws.Range("H22").Value = 1549.5555
$ws.Range("I22").Value = 1775
$ws.Range("J22").Value = 1369.2
$ws.Range("K22").Value = 5325
$ws.Range("L22").Value = 4107.6
$ws.Range("M22").Value = -5156
$ws.Range("N22").Value = -4445.6
$ws.Range("H27").Value = 1549.5555
$ws.Range("I27").Value = 1775
$ws.Range("J27").Value = 1369.2
$ws.Range("K27").Value = 5325
$ws.Range("L27").Value = 4107.6
$ws.Range("M27").Value = -5223
$ws.Range("N27").Value = -4311.6
$ws.Range("H34").Value = 2401.0715
$ws.Range("J34").Value = 2487.5
$ws.Range("L34").Value = 7462.5
$ws.Range("N34").Value = -7630.5
$ws.Range("H38").Value = 157.9
$ws.Range("J38").Value = 362
$ws.Range("L38").Value = 1086
$ws.Range("N38").Value = -1780
$ws.Range("H109").Value = 5397.579
$ws.Range("I109").Value = 2265.9092
$ws.Range("J109").Value = 9703.625
$ws.Range("K109").Value = 6797.7276
$ws.Range("L109").Value = 29110.875
$ws.Range("M109").Value = -5757.7276
$ws.Range("N109").Value = -31190.875
$ws.Range("H113").Value = 857.35297
$ws.Range("J113").Value = 894.1111
$ws.Range("L113").Value = 2682.3333
$ws.Range("N113").Value = -7022.3333
$ws.Range("H129").Value = 1042.381
$ws.Range("J129").Value = 2975
$ws.Range("L129").Value = 8925
$ws.Range("N129").Value = -18925
$ws.Range("H132").Value = 4946.125
$ws.Range("I132").Value = 1618.375
$ws.Range("J132").Value = 8273.875
$ws.Range("K132").Value = 14565.375
$ws.Range("L132").Value = 74464.875
$ws.Range("M132").Value = -12035.375
$ws.Range("N132").Value = -79524.875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 74749.5
$ws.Range("J63").Value = 74749.5
$ws.Range("L63").Value = 74749.5
$ws.Range("N63").Value = -76121.5
$ws.Range("H66").Value = 74749.5
$ws.Range("J66").Value = 74749.5
$ws.Range("L66").Value = 224248.5
$ws.Range("N66").Value = -231112.5
$ws.Range("H70").Value = 7461.407
$ws.Range("I70").Value = 7372.524
$ws.Range("K70").Value = 7372.524
$ws.Range("M70").Value = -7102.524
$ws.Range("H73").Value = 7461.407
$ws.Range("I73").Value = 7372.524
$ws.Range("K73").Value = 7372.524
$ws.Range("M73").Value = -6436.524
$ws.Range("H80").Value = 17858.285
$ws.Range("I80").Value = 8801.6
$ws.Range("J80").Value = 40500
$ws.Range("K80").Value = 8801.6
$ws.Range("L80").Value = 40500
$ws.Range("M80").Value = -7803.6
$ws.Range("N80").Value = -42496
$ws.Range("H83").Value = 17858.285
$ws.Range("I83").Value = 8801.6
$ws.Range("J83").Value = 40500
$ws.Range("K83").Value = 44008
$ws.Range("L83").Value = 202500
$ws.Range("M83").Value = -39016
$ws.Range("N83").Value = -212484
$ws.Range("H97").Value = 1081.8857
$ws.Range("I97").Value = 936.2
$ws.Range("J97").Value = 1276.1333
$ws.Range("K97").Value = 936.2
$ws.Range("L97").Value = 1276.1333
$ws.Range("M97").Value = -440.2
$ws.Range("N97").Value = -2268.1333
$ws.Range("H113").Value = 111113450
$ws.Range("I113").Value = 111113450
$ws.Range("K113").Value = 111113450
$ws.Range("M113").Value = -111111280
$ws.Range("H122").Value = 792274
$ws.Range("I122").Value = 2202799
$ws.Range("J122").Value = 8649
$ws.Range("K122").Value = 6608397
$ws.Range("L122").Value = 25947
$ws.Range("M122").Value = -6605947
$ws.Range("N122").Value = -30847
$ws.Range("H126").Value = 3904.2666
$ws.Range("I126").Value = 2408.611
$ws.Range("J126").Value = 6147.75
$ws.Range("K126").Value = 7225.833
$ws.Range("L126").Value = 18443.25
$ws.Range("M126").Value = -4755.833
$ws.Range("N126").Value = -23383.25
$ws.Range("H132").Value = 417484.06
$ws.Range("I132").Value = 80501.53999999999
$ws.Range("J132").Value = 3337999.2
$ws.Range("K132").Value = 241504.62
$ws.Range("L132").Value = 10013997.6
$ws.Range("M132").Value = -238974.62
$ws.Range("N132").Value = -10019057.6
$ws.Range("H133").Value = 124982
$ws.Range("J133").Value = 124982
$ws.Range("L133").Value = 124982
$ws.Range("N133").Value = -135102
$ws.Range("H135").Value = 85390
$ws.Range("J135").Value = 85390
$ws.Range("L135").Value = 85390
$ws.Range("N135").Value = -95530

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2854.48
$ws.Range("I16").Value = 2132.0588
$ws.Range("K16").Value = 2132.0588
$ws.Range("M16").Value = -1962.0588
$ws.Range("H22").Value = 1029.8
$ws.Range("I22").Value = 869.4
$ws.Range("J22").Value = 1511
$ws.Range("K22").Value = 869.4
$ws.Range("L22").Value = 1511
$ws.Range("M22").Value = -574.4
$ws.Range("N22").Value = -2101
$ws.Range("H27").Value = 1029.8
$ws.Range("I27").Value = 869.4
$ws.Range("J27").Value = 1511
$ws.Range("K27").Value = 869.4
$ws.Range("L27").Value = 1511
$ws.Range("M27").Value = -762.4
$ws.Range("N27").Value = -1725
$ws.Range("H40").Value = 4928.2354
$ws.Range("I40").Value = 2817.1
$ws.Range("K40").Value = 2817.1
$ws.Range("M40").Value = -2681.1
$ws.Range("H59").Value = 29950
$ws.Range("J59").Value = 29950
$ws.Range("L59").Value = 29950
$ws.Range("N59").Value = -31258
$ws.Range("H68").Value = 3501.2222
$ws.Range("I68").Value = 2752.0715
$ws.Range("J68").Value = 6123.25
$ws.Range("K68").Value = 2752.0715
$ws.Range("L68").Value = 6123.25
$ws.Range("M68").Value = -2003.0715
$ws.Range("N68").Value = -7621.25
$ws.Range("H71").Value = 3501.2222
$ws.Range("I71").Value = 2752.0715
$ws.Range("J71").Value = 6123.25
$ws.Range("K71").Value = 13760.3575
$ws.Range("L71").Value = 30616.25
$ws.Range("M71").Value = -10016.3575
$ws.Range("N71").Value = -38104.25
$ws.Range("H100").Value = 2112.25
$ws.Range("I100").Value = 1816.3334
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1816.3334
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1275.3334
$ws.Range("N100").Value = -4082
$ws.Range("H132").Value = 1999.38
$ws.Range("I132").Value = 1999.1171
$ws.Range("J132").Value = 2003.5
$ws.Range("K132").Value = 5997.3513
$ws.Range("L132").Value = 6010.5
$ws.Range("M132").Value = -3467.3513
$ws.Range("N132").Value = -11070.5
$ws.Range("H136").Value = 3137.3901
$ws.Range("I136").Value = 1925.75
$ws.Range("J136").Value = 7445.4443
$ws.Range("K136").Value = 5777.25
$ws.Range("L136").Value = 22336.3329
$ws.Range("M136").Value = -3227.25
$ws.Range("N136").Value = -27436.3329
$ws.Range("H137").Value = 31143
$ws.Range("J137").Value = 31143
$ws.Range("L137").Value = 31143
$ws.Range("N137").Value = -41343
$ws.Range("H139").Value = 78697.336
$ws.Range("J139").Value = 78697.336
$ws.Range("L139").Value = 78697.336
$ws.Range("N139").Value = -88977.336

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 6689.273
$ws.Range("I14").Value = 7426
$ws.Range("J14").Value = 5400
$ws.Range("K14").Value = 7426
$ws.Range("L14").Value = 5400
$ws.Range("M14").Value = -7258
$ws.Range("N14").Value = -5736
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30630
$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32184
$ws.Range("H96").Value = 3250
$ws.Range("I96").Value = 2500
$ws.Range("J96").Value = 4000
$ws.Range("K96").Value = 2500
$ws.Range("L96").Value = 4000
$ws.Range("M96").Value = -1127
$ws.Range("N96").Value = -6746
$ws.Range("H107").Value = 1222.5358
$ws.Range("I107").Value = 1409.5714
$ws.Range("K107").Value = 4228.7142
$ws.Range("M107").Value = -2308.7142
$ws.Range("H122").Value = 5033.4546
$ws.Range("I122").Value = 3670.0833
$ws.Range("J122").Value = 6669.5
$ws.Range("K122").Value = 11010.2499
$ws.Range("L122").Value = 20008.5
$ws.Range("M122").Value = -8560.249899999999
$ws.Range("N122").Value = -24908.5
$ws.Range("H126").Value = 4835.0713
$ws.Range("I126").Value = 4141
$ws.Range("K126").Value = 12423
$ws.Range("M126").Value = -9953
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 2867.1333
$ws.Range("I132").Value = 855.9091
$ws.Range("J132").Value = 8398
$ws.Range("K132").Value = 2567.7273
$ws.Range("L132").Value = 25194
$ws.Range("M132").Value = -37.72730000000001
$ws.Range("N132").Value = -30254
$ws.Range("H136").Value = 7737.091
$ws.Range("I136").Value = 1237.5807
$ws.Range("K136").Value = 3712.7421
$ws.Range("M136").Value = -1162.7421
